# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N on the "Repayment schedule"
# sheet - shifting the old "Late"/"Paid Date"/"Outstanding" columns one
# place to the right - and restores the cursor/selection state that was
# left behind on several sheets during the editing session.

$wb = $excel.ActiveWorkbook

# --- Repayment schedule: insert a blank column before N -------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Activate()

$wsSchedule.Columns("N").Insert()

# New column picks up the width of its left-hand neighbour (column M)
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

$wsSchedule.Range("R8").Select()

# --- Acc_Repayment: selection left on row 5 downward ----------------------
$wsAccRepayment = $wb.Worksheets.Item("Acc_Repayment")
$wsAccRepayment.Activate()
$wsAccRepayment.Range("A5:A1048576").EntireRow.Select()

# --- Acc_Disbursement1: selection left on D7 -------------------------------
$wsAccDisbursement1 = $wb.Worksheets.Item("Acc_Disbursement1")
$wsAccDisbursement1.Activate()
$wsAccDisbursement1.Range("D7").Select()

# --- Acc_Repayment1: selection left on row 5 downward ----------------------
$wsAccRepayment1 = $wb.Worksheets.Item("Acc_Repayment1")
$wsAccRepayment1.Activate()
$wsAccRepayment1.Range("A5:A1048576").EntireRow.Select()

# --- Acc_Upfront: stays the active sheet, selection left on F5 -------------
$wsAccUpfront = $wb.Worksheets.Item("Acc_Upfront")
$wsAccUpfront.Activate()
$wsAccUpfront.Range("F5").Select()
